$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove C2 and E2 entirely (no longer forecast values for this row)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: remove C3 entirely; tiny precision fix to E3
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = -0.6367039903685923

# Row 4: tiny precision fix to C4
$ws.Range("C4").Value = -3.956152295564885

# Row 5: tiny precision fixes to C5 and E5
$ws.Range("C5").Value = 1.234995474941436
$ws.Range("E5").Value = -0.209816187756795

# Row 6: tiny precision fix to C6
$ws.Range("C6").Value = 0.8993608108207818

# Row 8: tiny precision fix to C8
$ws.Range("C8").Value = 0.02019328874802717

# Row 10: tiny precision fix to E10
$ws.Range("E10").Value = -0.001769149545449711

# Row 12: tiny precision fix to C12
$ws.Range("C12").Value = 0.0720185131838802

# Row 13: tiny precision fix to E13
$ws.Range("E13").Value = -0.950584780912811

# Row 14: tiny precision fixes to C14 and E14
$ws.Range("C14").Value = -0.8017595264762423
$ws.Range("E14").Value = 0.0476740348578808

# Row 16: tiny precision fixes to C16 and E16
$ws.Range("C16").Value = 0.9704846793491706
$ws.Range("E16").Value = -0.8754609427830351

# Row 18: tiny precision fix to C18
$ws.Range("C18").Value = 0.3928252664241683

# Row 19: tiny precision fixes to C19 and E19
$ws.Range("C19").Value = 0.3224026462283369
$ws.Range("E19").Value = -0.9749878381046684
